$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.583.69"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "1.924.61"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.71"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4874"
$ws.Range("E7").Value = "  +2.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2911"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06733"
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "110.93"
$ws.Range("E10").Value = "  +4.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.15"
$ws.Range("E11").Value = "  +3.76%  "
$ws.Range("D12").Value = "1.919.62"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07591"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.315"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6720"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "296.70"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("D17").Value = "30.577.90"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.06"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007589"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.564"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").Value = "2.166.84"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9986"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.487"
$ws.Range("E24").Value = "  +2.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.484"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.69"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.30"
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.119"
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.456"
$ws.Range("E30").Value = "  +6.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.167"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.069"
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05041"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7422"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.142"
$ws.Range("E35").Value = "  -1.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9989"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02033"
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.686"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.59"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.029"
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4450"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8686"
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.03"
$ws.Range("E44").Value = "  +5.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.839"
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.278"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.66"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.263"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1236"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2547"
$ws.Range("E51").Value = "  +3.76%  "
